$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (bold font, border, centered
# alignment - style index 1) onto the two new header cells so they reuse
# the same style instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values (plain, unstyled like the rest of row 2)
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8
